# Add a new data row (row 99, columns B:L) to the "disinf" sheet, add column
# width customizations, and update the saved scroll/selection state — mirrors
# a single new record appended at the bottom of the log plus a new shared
# string ("3.45 PM") used by that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("disinf")

# --- Row 99: copy the date format (style) used by the row above first, ---
# --- then fill in all of the new row's values.                         ---
$ws.Range("G98").Copy()
$ws.Range("G99").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B99").Value = "Cheers"
$ws.Range("C99").Value = "Alwarpet"
$ws.Range("D99").Value = "Supermarket-Outlet"
$ws.Range("E99").Value = "Done "
$ws.Range("F99").Value = "Dinesh "
$ws.Range("G99").Value = 44603
$ws.Range("H99").Value = "3.45 PM"
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = "Self"
$ws.Range("K99").Value = "Fumigation "
$ws.Range("L99").Value = "CRM "

# --- Column widths for B, C, D, K, L ---
$ws.Columns.Item(2).ColumnWidth = 12.666666666666666   # B -> ~13.57
$ws.Columns.Item(3).ColumnWidth = 15.666666666666666   # C -> ~16.43
$ws.Columns.Item(4).ColumnWidth = 20.166666666666668   # D -> 21
$ws.Columns.Item(11).ColumnWidth = 13.666666666666666  # K -> ~14.43
$ws.Columns.Item(12).ColumnWidth = 9.666666666666666   # L -> ~10.43

# --- Scroll the view down and move the active selection to H100 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 85
$win.ScrollColumn = 1
$ws.Range("H100").Select()
